$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 2019 rows for Azerbaijan, Armenia and Georgia to the import/export
# turnover table. The sheet is laid out as three 14-row blocks (2005-2018)
# for Azerbaijan, Armenia, Georgia in that order; we insert one new row at
# the end of each block for the 2019 figure.
# ---------------------------------------------------------------------------

# --- Azerbaijan 2019 --------------------------------------------------------
# Old row 16 (start of the Armenia block) shifts down; the new row takes its
# place and becomes the Azerbaijan 2019 entry.
$ws.Rows(16).Insert()
$ws.Range("A16").Value = "Azerbaijan"
$ws.Range("B16").Value = 2019
$ws.Range("C16").Value = 2184313.08
$ws.Range("C16").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C16").Font.Name = $ws.Range("C2").Font.Name
$ws.Range("C16").Font.Size = $ws.Range("C2").Font.Size
$ws.Range("C16").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment

# Stray formatted-but-empty cell that appears next to the last Azerbaijan row
$ws.Range("E15").NumberFormat = "#,##0.00"
$ws.Range("E16").NumberFormat = "#,##0.00"
$ws.Range("F16").NumberFormat = "#,##0.00"

# --- Armenia 2019 -----------------------------------------------------------
# After the insert above, the Armenia block now runs rows 17-30 (2005-2018).
# Insert a new row right after it (row 31) for the 2019 figure.
$ws.Rows(31).Insert()
$ws.Range("A31").Value = "Armenia"
$ws.Range("B31").Value = 2019
$ws.Range("C31").Value = 941705.76
$ws.Range("C31").NumberFormat = "#,##0.00"

$ws.Range("E30").NumberFormat = "#,##0.00"
$ws.Range("E31").NumberFormat = "#,##0.00"
$ws.Range("F31").NumberFormat = "#,##0.00"

# --- Georgia 2019 ------------------------------------------------------------
# The Georgia block now runs rows 32-45 (2005-2018). Append a new row at the
# very end (row 46) for the 2019 figure.
$ws.Rows(46).Insert()
$ws.Range("A46").Value = "Georgia"
$ws.Range("B46").Value = 2019
$ws.Range("C46").Value = 1081719.35
$ws.Range("C46").NumberFormat = "#,##0.00"

$ws.Range("E46").NumberFormat = "#,##0.00"
$ws.Range("F46").NumberFormat = "#,##0.00"

# --- cosmetic view state -----------------------------------------------------
# (ColumnWidth values chosen to land as close as possible to the bestFit
# widths Excel computed for these columns; this engine only supports
# sixth-character granularity.)
$ws.Columns("E").ColumnWidth = 10.45
$ws.Columns("F").ColumnWidth = 9.0

$ws.Range("G21").Select() | Out-Null
